$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in the response data (mirrors row 2 for F/G/H/I) ---
$ws.Range("F3").Value = "123456789"
$ws.Range("G3").Value = "prueba"
$ws.Range("H3").Value = 1234
$ws.Range("I3").Value = 1234

# K3 needs a brand new shared string ("¡Lo Sentimos!") while preserving its
# existing cell style (which uses a quotePrefix number format). Writing the
# value directly via .Value/.Formula causes the engine to silently drop the
# quotePrefix flag from the style, creating an extra unwanted style record.
# To avoid that, stage the text in a scratch cell (default style), copy just
# the VALUE into K3 (leaving K3's own formatting untouched), then clean up
# the scratch cell.
$scratch = $ws.Range("Z100")
$scratch.Value = "¡Lo Sentimos!"
$scratch.Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$scratch.ClearContents() | Out-Null
$excel.CutCopyMode = 0

# --- New row 4: empty K4 cell with its own (underlined) font style ---
$ws.Range("K4").Value = ""
$ws.Range("K4").Font.Underline = 2                # xlUnderlineStyleSingle

# --- Sheet view: move the active selection to F3:I3 ---
$ws.Activate()
$ws.Range("F3:I3").Select() | Out-Null

Write-Host "edit applied"
